$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2026-01-27", "글로벌(OpenAI)", "How Indeed uses AI to help evolve the job search", "Mon, 26 Jan 2026 00:00:00 GMT", "https://openai.com/index/indeed-maggie-hulce"),
    @("2026-01-27", "글로벌(OpenAI)", "Unrolling the Codex agent loop", "Fri, 23 Jan 2026 12:00:00 GMT", "https://openai.com/index/unrolling-the-codex-agent-loop"),
    @("2026-01-27", "글로벌(OpenAI)", "Scaling PostgreSQL to power 800 million ChatGPT users", "Thu, 22 Jan 2026 12:00:00 GMT", "https://openai.com/index/scaling-postgresql"),
    @("2026-01-27", "글로벌(OpenAI)", "Inside Praktika's conversational approach to language learning", "Thu, 22 Jan 2026 05:00:00 GMT", "https://openai.com/index/praktika"),
    @("2026-01-27", "글로벌(OpenAI)", "Inside GPT-5 for Work: How Businesses Use GPT-5", "Thu, 22 Jan 2026 00:00:00 GMT", "https://openai.com/business/guides-and-resources/chatgpt-usage-and-adoption-patterns-at-work")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $cell = $ws.Cells.Item($row, $j + 1)
        $text = $values[$j]
        if ($j -eq 0) {
            # Column A holds plain "yyyy-mm-dd" strings (e.g. "2026-01-27").
            # Excel's normal value-entry heuristics would silently convert
            # that literal text into a real date serial number. Briefly
            # force a text format so the literal string is kept, then drop
            # the temporary formatting again so the cell ends up with no
            # explicit style, just like every other plain cell here.
            $cell.NumberFormatLocal = "@"
            $cell.Value = $text
            $cell.ClearFormats()
        } else {
            $cell.Value = $text
        }
    }
}
